# LoadFlow tang centering (not working though)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row labels in column A (adds two new shared strings: VrQ, Kg) ---
$ws.Range("A8").Value = "VrQ"
$ws.Range("A9").Value = "Kg"

# --- Updated input values ---
$ws.Range("B2").Value = 0
$ws.Range("B4").Value = 50

# --- New input supporting the Kg / VrQ centering logic ---
$ws.Range("B9").Value = 0.000001

# --- New helper formulas ---
$ws.Range("B8").Formula = "=MIN(B3 - B9,MAX(0,B2 + B9))"
$ws.Range("B10").Formula = "=(B8-B6)/B5"
$ws.Range("C2").Formula = "=B2*(1-B9)"
$ws.Range("C3").Formula = "=B3*(1-B9)"

# --- Updated offs formula (B7), now centers around the Kg-adjusted VrQ point ---
$ws.Range("B7").Formula = "=-ATANH((B8-B6)/B5)/B4"

# Force a full recalculation so every dependent cell (E column etc.) picks
# up the new inputs/formulas.
$excel.CalculateFull()

# --- Move active selection to B3 ---
$ws.Range("B3").Select()
